# Updated cryptos list (Price + Volume(1h) columns) on Tue Aug 13 02:59:59 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.249.85"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "2.679.08"
$ws.Range("E3").Value = "  +5.48%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").NumberFormat = "General"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'519.19"
$ws.Range("D5").NumberFormat = "General"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").Value = "'145.72"
$ws.Range("D6").NumberFormat = "General"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("D7").Value = "'0.994"
$ws.Range("D7").NumberFormat = "General"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "'0.570"
$ws.Range("D8").NumberFormat = "General"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").Value = "2.713.99"
$ws.Range("E9").Value = "  +6.78%  "
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").Value = "3.150.15"
$ws.Range("E14").Value = "  +5.72%  "
$ws.Range("D15").Value = "59.193.29"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "'21.17"
$ws.Range("D16").NumberFormat = "General"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "2.712.05"
$ws.Range("E18").Value = "  +6.83%  "
$ws.Range("D19").Value = "'356.53"
$ws.Range("D19").NumberFormat = "General"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  +6.28%  "
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "'10.48"
$ws.Range("D21").NumberFormat = "General"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "  +3.96%  "
$ws.Range("D22").Value = "'6.24"
$ws.Range("D22").NumberFormat = "General"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "  +5.02%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'62.10"
$ws.Range("D24").NumberFormat = "General"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("E25").Value = "  +3.69%  "
$ws.Range("D26").Value = "'0.991"
$ws.Range("D26").NumberFormat = "General"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("E28").Value = "  +3.56%  "
$ws.Range("E29").Value = "  +4.40%  "
$ws.Range("D30").Value = "'0.997"
$ws.Range("D30").NumberFormat = "General"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "'6.39"
$ws.Range("D31").NumberFormat = "General"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "  +9.42%  "
$ws.Range("D32").Value = "'19.16"
$ws.Range("D32").NumberFormat = "General"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  +3.47%  "
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("D34").Value = "'150.82"
$ws.Range("D34").NumberFormat = "General"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'0.977"
$ws.Range("D35").NumberFormat = "General"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("E36").Value = "  +3.13%  "
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("D39").Value = "'0.851"
$ws.Range("D39").NumberFormat = "General"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "  +3.08%  "
$ws.Range("E40").Value = "  +6.38%  "
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("D42").Value = "'283.96"
$ws.Range("D42").NumberFormat = "General"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").Value = "'0.621"
$ws.Range("D43").NumberFormat = "General"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("D44").Value = "'0.0992"
$ws.Range("D44").NumberFormat = "General"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "'19.95"
$ws.Range("D45").NumberFormat = "General"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("D46").Value = "'0.993"
$ws.Range("D46").NumberFormat = "General"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("D49").Value = "2.017.63"
$ws.Range("E49").Value = "  +6.64%  "
$ws.Range("E50").Value = "  +4.43%  "
$ws.Range("D51").Value = "'10.29"
$ws.Range("D51").NumberFormat = "General"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
